$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [object]$Range,
        [string]$Text
    )
    $escaped = $Text.Replace('"', '""')
    $Range.Formula = '="' + $escaped + '"'
    $Range.Copy()
    $Range.PasteSpecial(-4163)
}

function Set-CellText {
    param(
        [object]$Worksheet,
        [string]$Address,
        [string]$Text
    )
    Set-TextValue $Worksheet.Range($Address) $Text
}

# New ShipmentTracking numbers (column P) for rows 2-19, in row order
$trackingNumbers2to19 = @(
    @(2,  "320018586090"),
    @(3,  "320018586104"),
    @(4,  "320018586137"),
    @(5,  "320018586159"),
    @(6,  "320018586192"),
    @(7,  "320018586218"),
    @(8,  "320018586240"),
    @(9,  "320018586262"),
    @(10, "320018586295"),
    @(11, "320018586310"),
    @(12, "320018586354"),
    @(13, "320018586376"),
    @(14, "320018586402"),
    @(15, "320018586424"),
    @(16, "320018586457"),
    @(17, "320018586479"),
    @(18, "320018586516"),
    @(19, "320018586538")
)

foreach ($pair in $trackingNumbers2to19) {
    Set-CellText $ws "P$($pair[0])" $pair[1]
}

# Row 20 changed ActualRate (Q20), ShipmentTracking (P20) and Result (R20).
# The ActualRate was updated before the tracking number in the original edit.
Set-CellText $ws "Q20" "`$104.69"
Set-CellText $ws "P20" "320018586560"
$ws.Range("R20").Value = "FAIL"

# Remaining ShipmentTracking numbers (column P) for rows 21-26, in row order
$trackingNumbers21to26 = @(
    @(21, "320018586582"),
    @(22, "320018586619"),
    @(23, "320018586620"),
    @(24, "320018586630"),
    @(25, "320018586641"),
    @(26, "320018586652")
)

foreach ($pair in $trackingNumbers21to26) {
    Set-CellText $ws "P$($pair[0])" $pair[1]
}
